$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source row (row 3) holds the same values that need to be replicated into
# the new rows 4, 5 and 6.
$srcValues = $ws.Range("A3:AF3").Value2

foreach ($r in 4..6) {
    $destRange = $ws.Range("A" + $r + ":AF" + $r)
    $destRange.Value2 = $srcValues
}
